# Apply updated dSF (column F) values per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -4
    5  = -3
    8  = -3
    9  = -2
    14 = -3
    15 = -1
    16 = -3
    17 = -1
    18 = 4
    19 = -3
    20 = 3
    21 = -1
    22 = -1
    24 = 1
    25 = -2
    26 = -6
    27 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
